# This script applies the crypto price/volume data refresh described in the
# commit "Updated cryptos list ... with GitHub Actions".
#
# For each affected row it rewrites the Price (column D) and/or the
# Volume(1h) (column E) text values. Rows 49 and 50 additionally have their
# Coin name (B) and Link (C) values swapped, because "Quant" moved ahead of
# "PaxDollar" in the ranking.
#
# All of these source values are plain text (e.g. "0.998", "226.66") that
# Excel would otherwise auto-convert to numbers when assigned directly via
# .Value (losing significant trailing zeros / changing the cell type from
# inline string to numeric). To preserve the original text semantics we
# temporarily force the cell's number format to Text ("@") before writing
# the value, then clear the formatting again so no stray style id is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='34.802.08'}
    @{Cell='E2'; Value='  +3.07%  '}
    @{Cell='D3'; Value='1.805.76'}
    @{Cell='E3'; Value='  +1.23%  '}
    @{Cell='D4'; Value='0.998'}
    @{Cell='E4'; Value='  -0.44%  '}
    @{Cell='D5'; Value='226.66'}
    @{Cell='E5'; Value='  +0.98%  '}
    @{Cell='E6'; Value='  +1.00%  '}
    @{Cell='D7'; Value='0.997'}
    @{Cell='E7'; Value='  -0.47%  '}
    @{Cell='D8'; Value='33.37'}
    @{Cell='E8'; Value='  +7.47%  '}
    @{Cell='D9'; Value='0.287'}
    @{Cell='E9'; Value='  +3.09%  '}
    @{Cell='E10'; Value='  +2.92%  '}
    @{Cell='E11'; Value='  +1.43%  '}
    @{Cell='D12'; Value='2.063.70'}
    @{Cell='E12'; Value='  +1.19%  '}
    @{Cell='D13'; Value='11.37'}
    @{Cell='E13'; Value='  +14.69%  '}
    @{Cell='D14'; Value='1.824.43'}
    @{Cell='E14'; Value='  +2.28%  '}
    @{Cell='D15'; Value='0.643'}
    @{Cell='E15'; Value='  +2.68%  '}
    @{Cell='D16'; Value='34.818.96'}
    @{Cell='E16'; Value='  +3.12%  '}
    @{Cell='E17'; Value='  +3.74%  '}
    @{Cell='D18'; Value='69.84'}
    @{Cell='E18'; Value='  +2.04%  '}
    @{Cell='D19'; Value='257.67'}
    @{Cell='E19'; Value='  +2.91%  '}
    @{Cell='D20'; Value='0.0₃0774'}
    @{Cell='E20'; Value='  +5.32%  '}
    @{Cell='D21'; Value='0.999'}
    @{Cell='E21'; Value='  -0.27%  '}
    @{Cell='D22'; Value='10.56'}
    @{Cell='E22'; Value='  +3.13%  '}
    @{Cell='D23'; Value='4.27'}
    @{Cell='E23'; Value='  +1.66%  '}
    @{Cell='E24'; Value='  -0.36%  '}
    @{Cell='D25'; Value='158.99'}
    @{Cell='E25'; Value='  +0.57%  '}
    @{Cell='D26'; Value='16.59'}
    @{Cell='E26'; Value='  +1.46%  '}
    @{Cell='D27'; Value='7.18'}
    @{Cell='E27'; Value='  +4.27%  '}
    @{Cell='E28'; Value='  +0.24%  '}
    @{Cell='D29'; Value='0.998'}
    @{Cell='E29'; Value='  -0.50%  '}
    @{Cell='D30'; Value='3.83'}
    @{Cell='E30'; Value='  +1.26%  '}
    @{Cell='D31'; Value='0.0524'}
    @{Cell='E31'; Value='  +2.96%  '}
    @{Cell='E32'; Value='  +0.83%  '}
    @{Cell='E33'; Value='  +3.48%  '}
    @{Cell='D34'; Value='1.92'}
    @{Cell='E34'; Value='  +10.98%  '}
    @{Cell='D35'; Value='1.466.85'}
    @{Cell='E35'; Value='  -0.91%  '}
    @{Cell='E36'; Value='  +0.86%  '}
    @{Cell='D37'; Value='0.640'}
    @{Cell='E37'; Value='  +3.50%  '}
    @{Cell='D38'; Value='0.0191'}
    @{Cell='E38'; Value='  +3.06%  '}
    @{Cell='D39'; Value='84.42'}
    @{Cell='E39'; Value='  +2.59%  '}
    @{Cell='E40'; Value='  +3.73%  '}
    @{Cell='E41'; Value='  -1.22%  '}
    @{Cell='D42'; Value='0.913'}
    @{Cell='E42'; Value='  +3.83%  '}
    @{Cell='E43'; Value='  +2.23%  '}
    @{Cell='D44'; Value='0.0508'}
    @{Cell='E44'; Value='  -0.83%  '}
    @{Cell='D45'; Value='6.01'}
    @{Cell='E45'; Value='  +5.76%  '}
    @{Cell='D46'; Value='1.964.64'}
    @{Cell='E46'; Value='  +1.61%  '}
    @{Cell='E47'; Value='  -2.56%  '}
    @{Cell='D48'; Value='12.09'}
    @{Cell='E48'; Value='  +2.20%  '}
    @{Cell='B49'; Value='Quant'}
    @{Cell='C49'; Value='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'}
    @{Cell='D49'; Value='103.65'}
    @{Cell='E49'; Value='  +5.43%  '}
    @{Cell='B50'; Value='PaxDollar'}
    @{Cell='C50'; Value='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'}
    @{Cell='D50'; Value='0.997'}
    @{Cell='E50'; Value='  -0.52%  '}
    @{Cell='D51'; Value='50.48'}
    @{Cell='E51'; Value='  +0.05%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
